# TC02_Verify_MYACC_RegisteredUser.xlsx - "DOM Changes in ECTEST"
#
# Renames a few "Object" / "DataObject" test-data values on the two sheets
# (the QA author tweaked the element names used by the automation DOM
# lookups) and leaves the workbook with the "Testdata" sheet active/selected,
# mirroring how it was saved.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("TC02_Verify_MYACC_RegisteredUse")
$wsData = $wb.Worksheets.Item("Testdata")

# --- Update the Object column on the main testcase sheet ---------------
$wsMain.Range("C18").Value = "PurchasingHistory"
$wsMain.Range("C19").Value = "QuickOrderMyacc"
$wsMain.Range("C20").Value = "Storerooms"

# --- Update the matching DataObject values on the Testdata sheet -------
$wsData.Range("B14").Value = "Purchasing History"
$wsData.Range("B16").Value = "Storerooms"

# --- Restore the view/selection state saved with the workbook ----------
# Main sheet: scrolled down a bit, with the selection left on E23.
$wsMain.Activate() | Out-Null
$wsMain.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$wsMain.Range("E23").Select() | Out-Null

# Testdata sheet ends up the active/selected tab, with A18 selected.
$wsData.Activate() | Out-Null
$wsData.Range("A18").Select() | Out-Null

Write-Output "Updated DOM object names on $($wsMain.Name) and $($wsData.Name)."
